# Update "paises.xlsx" country stats + refresh timestamp.
# The source table (sheet "Pais") is sorted descending by "Casos totales" (col B).
# A data refresh bumped several countries' numbers which, in a few spots,
# changed the sort order locally (Aruba / Martinica / Curazao each moved up
# one row past their neighbour). We therefore rewrite each affected row's
# country name (col A) together with its full stats (cols B:H) so both the
# values and the row ordering end up correct.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Country, $Total, $New, $Active, $Recovered, $Critical, $DeathsToday, $Deaths) {
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $New
    $ws.Cells.Item($Row, 4).Value = $Active
    $ws.Cells.Item($Row, 5).Value = $Recovered
    $ws.Cells.Item($Row, 6).Value = $Critical
    $ws.Cells.Item($Row, 7).Value = $DeathsToday
    $ws.Cells.Item($Row, 8).Value = $Deaths
}

# Row 4 — Estados Unidos: refreshed totals
Set-Row 4 "Estados Unidos" 6513678 28008 3792099 2527565 0 480 194014

# Row 36 — Panama: refreshed totals
Set-Row 36 "Panama" 98407 829 70751 25549 0 8 2107

# Row 76 — Paraguay: refreshed totals
Set-Row 76 "Paraguay" 24214 861 11920 11831 0 14 463

# Row 115 — Congo: refreshed totals
Set-Row 115 "Congo" 4891 0 3887 890 0 12 114

# Rows 136-140 — Aruba jumps ahead of Bahamas/Jordania/Estonia/Sudan del Sur
Set-Row 136 "Aruba" 2589 107 1293 1281 0 0 15
Set-Row 137 "Bahamas" 2585 39 976 1550 0 1 59
Set-Row 138 "Jordania" 2581 103 1885 677 0 2 19
Set-Row 139 "Estonia" 2564 32 2195 305 0 0 64
Set-Row 140 "Sudan del Sur" 2552 7 1290 1213 0 1 49

# Rows 167-169 — Martinica jumps ahead of Santo Tome y Principe/Polinesia Francesa
Set-Row 167 "Martinica" 939 181 98 823 0 0 18
Set-Row 168 "Santo Tome y Principe" 898 0 859 24 0 0 15
Set-Row 169 "Polinesia Francesa" 773 0 528 245 0 0 0

# Row 189 — Barbados: refreshed totals
Set-Row 189 "Barbados" 180 1 156 17 0 0 7

# Row 191 — Monaco: refreshed totals
Set-Row 191 "Monaco" 156 3 107 48 0 0 1

# Rows 194-195 — Curazao jumps ahead of Liechtenstein
Set-Row 194 "Curazao" 107 5 48 58 0 0 1
Set-Row 195 "Liechtenstein" 107 0 97 9 0 0 1

# Refresh the "last updated" banner in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 02:55"
